$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has weekly price records ending at row 83 (dimension A1:R83).
# This week's update:
#  - inserts two new rows (so the sheet grows to A1:R85)
#  - adds a brand-new record at row 80 (Calameño / Super, fecha 44615)
#  - updates the existing "Extra" record's fecha + prices to the new week (row 81, fecha 44615)
#  - the remaining previously-existing records (Primera @ $/caja 16, Extra, Primera, Super)
#    shift down to rows 82-85 unchanged

# Insert two new rows right before the old row 80, pushing everything down.
$ws.Rows("80:81").Insert()

# New row 80: brand-new "Calameño / Super" record for fecha 44615
$ws.Range("A80").Value = 8
$ws.Range("B80").Value = "Terminal La Palmera de La Serena"
$ws.Range("C80").Value = "Coquimbo"
$ws.Range("D80").Value = 44615
$ws.Range("E80").Value = 4
$ws.Range("F80").Value = 100112027
$ws.Range("G80").Value = "Melón"
$ws.Range("H80").Value = "Calameño"
$ws.Range("I80").Value = "Super"
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 1700
$ws.Range("L80").Value = 1800
$ws.Range("M80").Value = 1750
$ws.Range("N80").Value = "$/unidad"
$ws.Range("O80").Value = "Región de O'Higgins"
$ws.Range("P80").Value = 1750
$ws.Range("Q80").Value = 1
$ws.Range("R80").Value = "Hortaliza"

# New row 81: "Tuna / Super" record, fecha updated to 44615, prices updated
$ws.Range("A81").Value = 8
$ws.Range("B81").Value = "Terminal La Palmera de La Serena"
$ws.Range("C81").Value = "Coquimbo"
$ws.Range("D81").Value = 44615
$ws.Range("E81").Value = 4
$ws.Range("F81").Value = 100112027
$ws.Range("G81").Value = "Melón"
$ws.Range("H81").Value = "Tuna"
$ws.Range("I81").Value = "Super"
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 1700
$ws.Range("L81").Value = 1800
$ws.Range("M81").Value = 1750
$ws.Range("N81").Value = "$/unidad"
$ws.Range("O81").Value = "Región de O'Higgins"
$ws.Range("P81").Value = 1750
$ws.Range("Q81").Value = 1
$ws.Range("R81").Value = "Hortaliza"
